$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# --- Row 11: add a "default_position" (G) value of "body" ---
$ws.Cells.Item(11, 7).Value = "body"

# --- New row 190: Astral Ring (spell) ---
$ws.Cells.Item(190, 2).Value  = "Astral Ring"
$ws.Cells.Item(190, 3).Value  = 1
$ws.Cells.Item(190, 4).Value  = "Hate Filled Dreams"
$ws.Cells.Item(190, 5).Value  = "spell-damage"
$ws.Cells.Item(190, 6).Value  = "Cause the enemy to have hate filled dreams."
$ws.Cells.Item(190, 8).Value  = 85
$ws.Cells.Item(190, 11).Value = 18000
$ws.Cells.Item(190, 21).Value = 1
$ws.Cells.Item(190, 22).Value = 36
$ws.Cells.Item(190, 23).Value = 70
$ws.Cells.Item(190, 24).Value = "spell"

# --- New row 191: Fighters Strength (armour / helmet) ---
$ws.Cells.Item(191, 1).Value  = "Fighters Strength"
$ws.Cells.Item(191, 3).Value  = 1
$ws.Cells.Item(191, 4).Value  = "Bone Chip Helmet"
$ws.Cells.Item(191, 5).Value  = "helmet"
$ws.Cells.Item(191, 6).Value  = "Made of the enemies bones, we chip away at their remains for the pieces they offer us."
$ws.Cells.Item(191, 7).Value  = "helmet"
$ws.Cells.Item(191, 10).Value = 30
$ws.Cells.Item(191, 11).Value = 945
$ws.Cells.Item(191, 15).Value = 0.06
$ws.Cells.Item(191, 16).Value = 0.06
$ws.Cells.Item(191, 17).Value = 0.06
$ws.Cells.Item(191, 18).Value = 0.06
$ws.Cells.Item(191, 19).Value = 0.06
$ws.Cells.Item(191, 21).Value = 1
$ws.Cells.Item(191, 22).Value = 24
$ws.Cells.Item(191, 23).Value = 40
$ws.Cells.Item(191, 24).Value = "armour"

# --- New row 192: Natures Balancing Bliss (armour / body) ---
$ws.Cells.Item(192, 1).Value  = "Natures Balancing Bliss"
$ws.Cells.Item(192, 3).Value  = 1
$ws.Cells.Item(192, 4).Value  = "Ripped Cloth"
$ws.Cells.Item(192, 5).Value  = "body"
$ws.Cells.Item(192, 6).Value  = "It's at least clothing, to say the least."
$ws.Cells.Item(192, 7).Value  = "body"
$ws.Cells.Item(192, 10).Value = 4
$ws.Cells.Item(192, 11).Value = 10
$ws.Cells.Item(192, 21).Value = 1
$ws.Cells.Item(192, 22).Value = 1
$ws.Cells.Item(192, 23).Value = 5
$ws.Cells.Item(192, 24).Value = "armour"

# --- New row 193: Spell Crafters Blood (weapon) ---
$ws.Cells.Item(193, 2).Value  = "Spell Crafters Blood"
$ws.Cells.Item(193, 3).Value  = 1
$ws.Cells.Item(193, 4).Value  = "Warriors Battle Axe"
$ws.Cells.Item(193, 5).Value  = "weapon"
$ws.Cells.Item(193, 6).Value  = "Made for a warrior. This battle axe will cut down all your enemies."
$ws.Cells.Item(193, 8).Value  = 80
$ws.Cells.Item(193, 11).Value = 750
$ws.Cells.Item(193, 15).Value = 0.18
$ws.Cells.Item(193, 16).Value = 0.18
$ws.Cells.Item(193, 17).Value = 0.18
$ws.Cells.Item(193, 18).Value = 0.18
$ws.Cells.Item(193, 19).Value = 0.18
$ws.Cells.Item(193, 21).Value = 1
$ws.Cells.Item(193, 22).Value = 18
$ws.Cells.Item(193, 23).Value = 36
$ws.Cells.Item(193, 24).Value = "weapon"

# --- New row 194: Weapon Crafter Spell (ring) ---
$ws.Cells.Item(194, 2).Value  = "Weapon Crafter Spell"
$ws.Cells.Item(194, 3).Value  = 1
$ws.Cells.Item(194, 4).Value  = "Glass Ring"
$ws.Cells.Item(194, 5).Value  = "ring"
$ws.Cells.Item(194, 6).Value  = "Made from the finest glass and enforced by a metal band."
$ws.Cells.Item(194, 8).Value  = 75
$ws.Cells.Item(194, 11).Value = 133000
$ws.Cells.Item(194, 12).Value = 0.25
$ws.Cells.Item(194, 13).Value = 0.25
$ws.Cells.Item(194, 14).Value = 0.25
$ws.Cells.Item(194, 21).Value = 1
$ws.Cells.Item(194, 22).Value = 50
$ws.Cells.Item(194, 23).Value = 100
$ws.Cells.Item(194, 24).Value = "ring"

# --- Column width updates for columns A and B (now hold longer item names) ---
$ws.Columns.Item(1).ColumnWidth = 27.25
$ws.Columns.Item(2).ColumnWidth = 23.75
